$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Update column D values (column F recalculates automatically via shared formula)
$ws.Range("D205").Value = 39.793691562464502
$ws.Range("D206").Value = 88.284824261832398
$ws.Range("D207").Value = 23.679460210882201
$ws.Range("D208").Value = 46.034813829787197
$ws.Range("D209").Value = 29.993174327587798
$ws.Range("D210").Value = 29.232213618983401
$ws.Range("D212").Value = 28.5602246963563
$ws.Range("D213").Value = 59.321030355944799
$ws.Range("D214").Value = 59.0453053240619
$ws.Range("D215").Value = 122.35614210717701
$ws.Range("D217").Value = 29.747237435146999
$ws.Range("D221").Value = 28.405353891373501
$ws.Range("D223").Value = 46.965378395248202
$ws.Range("D225").Value = 49.311700677614397
$ws.Range("D227").Value = 31.736438859868301
$ws.Range("D228").Value = 41.849537980867403
$ws.Range("D229").Value = 130.534834693878
$ws.Range("D231").Value = 23.691518078357898
$ws.Range("D232").Value = 55.077469074596401
$ws.Range("D233").Value = 42.328923245331801
$ws.Range("D234").Value = 63.5346849245675
$ws.Range("D237").Value = 102.288730602621
$ws.Range("D239").Value = 51.154247895412702
$ws.Range("D240").Value = 63.175570181865297
$ws.Range("D241").Value = 123.130531421385
$ws.Range("D242").Value = 24.535693850118399
$ws.Range("D243").Value = 32.441431309253403
$ws.Range("D246").Value = 102.368161072514
$ws.Range("D249").Value = 108.27233589868
$ws.Range("D250").Value = 101.556414255474
$ws.Range("D256").Value = 31.3121581123612
$ws.Range("D258").Value = 33.8391124773825
$ws.Range("D260").Value = 38.411854172615499
$ws.Range("D261").Value = 122.35614210717701
$ws.Range("D263").Value = 44.7217585319747
$ws.Range("D265").Value = 77.743450460332099
$ws.Range("D266").Value = 118.05711026462799
$ws.Range("D268").Value = 107.57585908084199
$ws.Range("D271").Value = 40.127692513604899
$ws.Range("D272").Value = 44.4750088980825

# Update the view selection / scroll position to match the saved workbook state
$win = $excel.ActiveWindow
$win.ScrollRow = 209
$win.ScrollColumn = 1
$ws.Range("G218").Select()

